# feat: add 2022-Q4 data
#
# The "2022-Q3" sheet's data is rolled forward into a new "2022-Q4" sheet
# (reusing the original sheet object / position), a duplicate of the
# original "2022-Q3" sheet is created right after it to preserve the old
# Q3 numbers, and the "总计" (summary) sheet gets a new row for Q2 while
# Q3/Q4 rows shift down.

$wb = $excel.ActiveWorkbook

# --- 1. Duplicate the "2022-Q3" sheet so the original quarter's numbers
#        survive under their own tab, placed right after the source sheet.
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsQ3.Copy($null, $wsQ3)
$wsQ3Copy = $wb.Worksheets.Item("2022-Q3 (2)")

# --- 2. Turn the original "2022-Q3" sheet into the new "2022-Q4" sheet
#        (it keeps its tab position, right after "总计") and rename the
#        copy back to "2022-Q3".
$wsQ3.Name = "2022-Q4"
$wsQ3Copy.Name = "2022-Q3"

$wsQ4 = $wsQ3

# --- 3. Update the Q4 sheet's figures with the new quarter's data. These
#        columns hold text-formatted numbers (e.g. "6.60"), so force text
#        formatting before assigning or Excel would coerce them to numbers
#        and lose the trailing zero.
$wsQ4NumRange = $wsQ4.Range("D2:G3")
$wsQ4NumRange.NumberFormat = "@"

$wsQ4.Range("D2").Value = "0.21"
$wsQ4.Range("E2").Value = "86.58"
$wsQ4.Range("F2").Value = "6.60"
$wsQ4.Range("G2").Value = "0.0139"

$wsQ4.Range("D3").Value = "0.16"
$wsQ4.Range("E3").Value = "86.58"
$wsQ4.Range("F3").Value = "6.60"
$wsQ4.Range("G3").Value = "0.0106"

# --- 4. Update the "总计" summary sheet: insert a new row for the
#        (older) 2022-Q2 entry and shift the quarter labels/values so the
#        sheet again reads Q4 / Q3 / Q2 top to bottom.
$wsTotal = $wb.Worksheets.Item("总计")

$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("D2").Value = 0.02

$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("D3").Value = 0.03

# Row 4 is new — copy row 3's formatting down (matches the A-column style
# used by the existing rows) before filling in the Q2 values.
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2022-Q2"
$wsTotal.Range("C4").Value = 2
$wsTotal.Range("D4").Value = 0.02

# --- 5. The sheet Copy() operation above shifted the active tab; restore
#        the original active sheet ("2022-Q2").
$wb.Worksheets.Item("2022-Q2").Activate()
